# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E50) previously listed period codes in
# descending order (2003 down to 1705). The data got refreshed and the
# periods are now listed in ascending order (1705 up to 2003), and the
# "Valor Mora" column (F16:F50) was updated with the corresponding new
# figures for each period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period list for E16:E50 (was descending 2003->1705)
$periods = @(
    "1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

# New "Valor Mora" figures aligned with the refreshed periods above
$valores = @(
    22624,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,29509,
    31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249
)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
